# fix best run MAG9
# The "best_desman_run" value for MAG9 (row 10, column B) was recorded as
# MAG9_7_9 but should have been MAG9_7_1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "MAG9_7_1"

# Leave the selection where the user ended up after editing the cell.
$ws.Range("B11").Select()
